# Add a couple of new user stories (US04, US05) to the
# "Product backlog Pithia" worksheet, rows 8 and 9, columns A-C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product backlog Pithia")
$ws.Activate()

# --- Row 8: US04 ---------------------------------------------------------
$ws.Range("A8").Value = 'US04'
$ws.Range("B8").Value = 'Ως καθηγητης,θελω να εχω δυνατοτητα να βαζω βαθμολογια στους μαθητες που παρακολουθουν το μαθημα/μαθηματα που διδασκω.'
$ws.Range("C8").Value = 'Οταν ο καθηγητης μπαινει στο συστημα,θα υπαρχει στο μενου επιλογη "Δηλωση Βαθμολογιας" το οποιο θα οδηγει τον καθηγητη στη φορμα δηλωσης βαθμολογιας φοιτητων.'

# --- Row 9: US05 ---------------------------------------------------------
$ws.Range("A9").Value = 'US05'
$ws.Range("B9").Value = 'Ως διαχειριστης,θελω να εχω δυνατοτητα να προσθετω και να αφαιρω καθηγητη/-ες απο το συστημα.'
$ws.Range("C9").Value = 'Οταν ο διαχειριστης μπαινει στο συστημα,θα υπαρχει στο μενου επιλογη "Προσθηκη/Αφαιρεση Καθηγητη" το οποιο θα τον οδηγει στο περιβαλλον διαχειρισης καθηγητων.'

# --- Reflect the updated view state (scrolled down a bit, new selection) -
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("C10").Select()
